$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 265, shifting the existing rows 265:281 down to 266:282.
$ws.Rows("265:265").Insert()

# Populate the newly inserted row 265 with the new weekly price observation.
$ws.Cells.Item(265, 1).Value = 3
$ws.Cells.Item(265, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(265, 3).Value = "Coquimbo"
$ws.Cells.Item(265, 4).Value = 44585
$ws.Cells.Item(265, 5).Value = 5
$ws.Cells.Item(265, 6).Value = 100114013
$ws.Cells.Item(265, 7).Value = "Zanahoria"
$ws.Cells.Item(265, 8).Value = "Sin especificar"
$ws.Cells.Item(265, 9).Value = "Primera"
$ws.Cells.Item(265, 10).Value = 160
$ws.Cells.Item(265, 11).Value = 8000
$ws.Cells.Item(265, 12).Value = 8000
$ws.Cells.Item(265, 13).Value = 8000
$ws.Cells.Item(265, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(265, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(265, 16).Value = 400
$ws.Cells.Item(265, 17).Value = 20
$ws.Cells.Item(265, 18).Value = "Hortaliza"

# Keep the date column formatted consistently with the rest of column D.
$ws.Cells.Item(265, 4).NumberFormat = $ws.Cells.Item(266, 4).NumberFormat
